$wb = $excel.ActiveWorkbook

# --- Update workbook metadata (fileVersion rupBuild / revisionPtr documentId are
# engine/Excel-build identifiers not reachable via the Excel object model; skipped) ---

# --- PO List sheet: update tracked counter/date cells ---
$poList = $wb.Worksheets.Item("PO List")

$poList.Range("N3").Value = 14
$poList.Range("K4").Value = 7
$poList.Range("L4").Value = 7
$poList.Range("N4").Value = 19
$poList.Range("O4").Value = 3
$poList.Range("P4").Value = 3
$poList.Range("Q4").Value = 44967
$poList.Range("R4").Value = 2
$poList.Range("S4").Value = 2
$poList.Range("T4").Value = 2
$poList.Range("U4").Value = 44967
$poList.Range("V4").Value = 9
$poList.Range("W4").Value = 9
$poList.Range("X4").Value = 44967
$poList.Range("Y4").Value = 1
$poList.Range("K5").Value = 3
$poList.Range("L5").Value = 3
$poList.Range("N5").Value = 22
$poList.Range("O5").Value = 6
$poList.Range("P5").Value = 6
$poList.Range("Q5").Value = 44988
$poList.Range("R5").Value = 1
$poList.Range("S5").Value = 5
$poList.Range("T5").Value = 5
$poList.Range("U5").Value = 44988
$poList.Range("N6").Value = 23
$poList.Range("N7").Value = 3
$poList.Range("R7").Value = 5
$poList.Range("N8").Value = 21
$poList.Range("R8").Value = 6
$poList.Range("N9").Value = 17
$poList.Range("R9").Value = 8
$poList.Range("Y9").Value = 6
$poList.Range("N10").Value = 13
$poList.Range("N11").Value = 11
$poList.Range("R11").Value = 4
$poList.Range("Y11").Value = 3
$poList.Range("N12").Value = 9
$poList.Range("Y12").Value = 7
$poList.Range("N13").Value = 8
$poList.Range("N14").Value = 20
$poList.Range("N16").Value = 15
$poList.Range("R16").Value = 3
$poList.Range("Y16").Value = 2
$poList.Range("K18").Value = 4
$poList.Range("L18").Value = 4
$poList.Range("M18").Value = 44984
$poList.Range("N18").Value = 1
$poList.Range("R18").Value = 10
$poList.Range("N19").Value = 6
$poList.Range("Y19").Value = 5
$poList.Range("N20").Value = 3
$poList.Range("N21").Value = 18
$poList.Range("F22").Value = 29
$poList.Range("N24").Value = 10
$poList.Range("N25").Value = 7
$poList.Range("R25").Value = 7
$poList.Range("N26").Value = 16
$poList.Range("N27").Value = 5
$poList.Range("R27").Value = 11
$poList.Range("N28").Value = 2
$poList.Range("Y28").Value = 4
$poList.Range("N29").Value = 12
$poList.Range("R29").Value = 9

# --- Sheet-view / active-tab changes ---
# Select the bottom-right pane cell J39 on "PO List" (was B3:Z29).
$poList.Activate()
$poList.Range("J39").Select()

# Make "Assignment" the active/selected sheet (was "PO List").
$assignment = $wb.Worksheets.Item("Assignment")
$assignment.Activate()
